$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
# The previous account-statement periods (rows 17-26) are replaced with the
# new periods (now listed in descending/most-recent-first order) and their
# updated "Salario Basico" (F) / "Valor Mora" (G) amounts.
# ---------------------------------------------------------------------------

$data = @(
    @{Row=17; Periodo="2203"; Salario=297633; Mora=8267587},
    @{Row=18; Periodo="2202"; Salario=330703; Mora=8267587},
    @{Row=19; Periodo="2201"; Salario=330703; Mora=8267587},
    @{Row=20; Periodo="2112"; Salario=330703; Mora=8267587},
    @{Row=21; Periodo="2111"; Salario=330703; Mora=8267587},
    @{Row=22; Periodo="2110"; Salario=330703; Mora=8267587},
    @{Row=23; Periodo="2109"; Salario=330703; Mora=8267587},
    @{Row=24; Periodo="2108"; Salario=330703; Mora=8267587},
    @{Row=25; Periodo="2106"; Salario=330703; Mora=8267587},
    @{Row=26; Periodo="2105"; Salario=330703; Mora=8267587}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Range("E" + $r).Value = $item.Periodo
    $ws.Range("F" + $r).Value = $item.Salario
    $ws.Range("G" + $r).Value = $item.Mora
}

# ---------------------------------------------------------------------------
# Column widths were adjusted (wider) to accommodate the refreshed figures.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 17.583635416666667
$ws.Columns.Item(3).ColumnWidth = 15.750229166666667
$ws.Columns.Item(5).ColumnWidth = 12.583635416666667
$ws.Columns.Item(7).ColumnWidth = 13.416947916666667
$ws.Columns.Item(8).ColumnWidth = 18.41694791666667
$ws.Columns.Item(9).ColumnWidth = 17.250510416666668
$ws.Columns.Item(10).ColumnWidth = 14.083666666666666
